$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (TC_01): Priority HIGH, new Pre-conditions text ---
$ws.Range("B2").Value = "HIGH"
$ws.Range("E2").Value = "1. Aplikacja jest uruchomiona. `n 2. Dostępna jest lista produktów (menu) zawierająca produkty bez ograniczeń wiekowych oraz produkty oznaczone jako (18+)."

# --- Row 3 (TC_02): Priority HIGH, new Pre-conditions text ---
$ws.Range("B3").Value = "HIGH"
$ws.Range("E3").Value = "1. Aplikacja jest uruchomiona. `n 2. Baza danych produktów zawiera co najmniej jeden produkt z flaga czy_tylko_dla_pelnoletnich=True."

# --- Row 4 (TC_03): new Pre-conditions text ---
$ws.Range("E4").Value = "1. Aplikacja jest uruchomiona i wyświetla ekran startowy (rejestracja klienta). `n 2. System ma dostęp do aktualnej daty systemowej w celu wyliczenia wieku."
$ws.Rows.Item(4).RowHeight = 100.8

# --- Row 5 (TC_04): new Pre-conditions text ---
$ws.Range("E5").Value = "1. Aplikacja jest uruchomiona. `n 2. Klient został pomyślnie zarejestrowany i znajduje się na etapie wyboru produktów. `n 3. Tester zna zakres dostępnych ID w menu."

# --- Row 6 (new TC_05) ---
$ws.Range("A6").Value = "TC_05"
$ws.Range("B6").Value = "LOW"
$ws.Range("C6").Value = "Negative / Functional"
$ws.Range("D6").Value = "Weryfikacja błędnie wprowadzonych danych nie będących liczbami naturalnymi(znakiem z klawiatury)"
$ws.Range("E6").Value = "1. Klient jest zalogowany/zarejestrowany w systemie. `n 2. Pole wejściowe ID produktu oczekuje na wprowadzenie danych przez użytkownika."
$ws.Range("F6").Value = "1. Wyświetl listę dostępnych produktów (menu). `n 2. Zidentyfikuj najwyższe dostępne ID (np. 7). `n 3. W polu wyboru ID wpisz wartość o jeden znak nie będący liczbą naturalną."
$ws.Range("G6").Value = "System wyświetla komunikat ostrzegawczy i nie dodaje nic do koszyka."
$ws.Range("H6").Value = "System wyświetla komunikat ostrzegawczy i nie dodaje nic do koszyka."
$ws.Range("I6").Value = "PASS"
$ws.Range("A6:I6").HorizontalAlignment = -4108
$ws.Range("A6:I6").VerticalAlignment = -4108
$ws.Range("A6:I6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 72

# --- Selection / view state ---
$ws.Range("E5").Select()
